$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "完成情况" (completion status) column with "完成" for the
# five task rows that previously had it blank.
$ws.Range("C3").Value = "完成"
$ws.Range("C4").Value = "完成"
$ws.Range("C5").Value = "完成"
$ws.Range("C6").Value = "完成"
$ws.Range("C7").Value = "完成"

# Copy the cell formatting from the adjacent "计划内容" cell so the new
# text matches the existing table styling (bordered cell).
$ws.Range("B3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null

# The second block's "总结：" notes (previously listing admin-function
# bullet points) are cleared out, leaving the rows blank.
$ws.Range("B19").ClearContents()
$ws.Range("B20").ClearContents()

# Update the active selection to match where editing ended up.
$ws.Range("C7").Select()
